$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Trening" header in column F, matching the existing header style
$ws.Cells.Item(1, 6).Value = "Trening"
$ws.Range("A1").Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)


# Fill the data table (rows 2-13)

$ws.Cells.Item(2, 1).Value = 45684.59172523148
$ws.Cells.Item(2, 2).Value = 525
$ws.Cells.Item(2, 3).Value = 11.25
$ws.Cells.Item(2, 4).Value = 1.736757789339338
$ws.Cells.Item(2, 5).Value = "10-15"
$ws.Cells.Item(2, 6).Value = "Duża Gra"

$ws.Cells.Item(3, 1).Value = 45684.59349143518
$ws.Cells.Item(3, 2).Value = 677.6
$ws.Cells.Item(3, 3).Value = 11.49
$ws.Cells.Item(3, 4).Value = 1.993386421884807
$ws.Cells.Item(3, 5).Value = "10-15"
$ws.Cells.Item(3, 6).Value = "Duża Gra"

$ws.Cells.Item(4, 1).Value = 45684.59414537037
$ws.Cells.Item(4, 2).Value = 734.1
$ws.Cells.Item(4, 3).Value = 11.44
$ws.Cells.Item(4, 4).Value = 2.19777490411486
$ws.Cells.Item(4, 5).Value = "10-15"
$ws.Cells.Item(4, 6).Value = "Duża Gra"

$ws.Cells.Item(5, 1).Value = 45684.59210833333
$ws.Cells.Item(5, 2).Value = 558.1
$ws.Cells.Item(5, 3).Value = 9
$ws.Cells.Item(5, 4).Value = 1.734153338841029
$ws.Cells.Item(5, 5).Value = "5-10"
$ws.Cells.Item(5, 6).Value = "Duża Gra"

$ws.Cells.Item(6, 1).Value = 45684.59348912037
$ws.Cells.Item(6, 2).Value = 677.4
$ws.Cells.Item(6, 3).Value = 9.85
$ws.Cells.Item(6, 4).Value = 1.910983051572527
$ws.Cells.Item(6, 5).Value = "5-10"
$ws.Cells.Item(6, 6).Value = "Duża Gra"

$ws.Cells.Item(7, 1).Value = 45684.59414305555
$ws.Cells.Item(7, 2).Value = 733.9
$ws.Cells.Item(7, 3).Value = 9.74
$ws.Cells.Item(7, 4).Value = 2.041442751884461
$ws.Cells.Item(7, 5).Value = "5-10"
$ws.Cells.Item(7, 6).Value = "Duża Gra"

$ws.Cells.Item(8, 1).Value = 45684.59857013889
$ws.Cells.Item(8, 2).Value = 1116.4
$ws.Cells.Item(8, 3).Value = 10.29
$ws.Cells.Item(8, 4).Value = 2.801578249250139
$ws.Cells.Item(8, 5).Value = "10-15"
$ws.Cells.Item(8, 6).Value = "Mała Gra"

$ws.Cells.Item(9, 1).Value = 45684.60083518519
$ws.Cells.Item(9, 2).Value = 1312.1
$ws.Cells.Item(9, 3).Value = 11.19
$ws.Cells.Item(9, 4).Value = 2.749141931533814
$ws.Cells.Item(9, 5).Value = "10-15"
$ws.Cells.Item(9, 6).Value = "Mała Gra"

$ws.Cells.Item(10, 1).Value = 45684.60317893518
$ws.Cells.Item(10, 2).Value = 1514.6
$ws.Cells.Item(10, 3).Value = 10.88
$ws.Cells.Item(10, 4).Value = 3.320722034999303
$ws.Cells.Item(10, 5).Value = "10-15"
$ws.Cells.Item(10, 6).Value = "Mała Gra"

$ws.Cells.Item(11, 1).Value = 45684.59856898148
$ws.Cells.Item(11, 2).Value = 1116.3
$ws.Cells.Item(11, 3).Value = 9.42
$ws.Cells.Item(11, 4).Value = 2.894017015184676
$ws.Cells.Item(11, 5).Value = "5-10"
$ws.Cells.Item(11, 6).Value = "Mała Gra"

$ws.Cells.Item(12, 1).Value = 45684.60317777778
$ws.Cells.Item(12, 2).Value = 1514.5
$ws.Cells.Item(12, 3).Value = 9.54
$ws.Cells.Item(12, 4).Value = 3.221963201250348
$ws.Cells.Item(12, 5).Value = "5-10"
$ws.Cells.Item(12, 6).Value = "Mała Gra"

$ws.Cells.Item(13, 1).Value = 45684.6039949074
$ws.Cells.Item(13, 2).Value = 1585.1
$ws.Cells.Item(13, 3).Value = 9.55
$ws.Cells.Item(13, 4).Value = 2.803569725581578
$ws.Cells.Item(13, 5).Value = "5-10"
$ws.Cells.Item(13, 6).Value = "Mała Gra"

# Apply the date/time number format to column A data cells (2-13).
# First cell goes through the lowercase-then-uppercase format transition
# (mirrors the authoring session, leaving numFmtId 164 defined but unused),
# all subsequent cells reuse the resulting uppercase style directly.
$ws.Cells.Item(2, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(2, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
for ($r = 3; $r -le 13; $r++) {
  $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
